$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B class
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9466666666666667
$ws.Range("C2").Value = 0.9861111111111112
$ws.Range("D2").Value = 0.9659863945578231
$ws.Range("E2").Value = 72

# Row 3: M class
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9743589743589743
$ws.Range("C3").Value = 0.9047619047619048
$ws.Range("D3").Value = 0.9382716049382716
$ws.Range("E3").Value = 42

# Row 4: accuracy
$ws.Range("B4").Value = 0.956140350877193
$ws.Range("C4").Value = 0.956140350877193
$ws.Range("D4").Value = 0.956140350877193
$ws.Range("E4").Value = 0.956140350877193

# Row 5: macro avg
$ws.Range("B5").Value = 0.9605128205128205
$ws.Range("C5").Value = 0.9454365079365079
$ws.Range("D5").Value = 0.9521289997480473
$ws.Range("E5").Value = 114

# Row 6: weighted avg
$ws.Range("B6").Value = 0.9568690958164642
$ws.Range("C6").Value = 0.956140350877193
$ws.Range("D6").Value = 0.9557756825927252
$ws.Range("E6").Value = 114
